# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2304
#   *_new  -> *_FV2310
# Also freeze the header row and wrap the sheet's data range in an Excel
# Table (Table1) so the new headers double as the table's column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2304 = "_FV2304"
$fv2310 = "_FV2310"

$lastCol = $ws.UsedRange.Columns.Count
$lastRow = $ws.UsedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -ne $null) {
        if ($v.EndsWith($oldSuffix)) {
            $cell.Value = $v.Substring(0, $v.Length - $oldSuffix.Length) + $fv2304
        } elseif ($v.EndsWith($newSuffix)) {
            $cell.Value = $v.Substring(0, $v.Length - $newSuffix.Length) + $fv2310
        }
    }
}

# Freeze the header row (row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into an Excel table, reusing the freshly renamed
# header row as the table's column headers.
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
